$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-create the hyperlinks on E2/E3, keeping the same target mailto address
# and the original e-mail text as the hyperlink's display text. This must
# happen BEFORE the cell values are updated below, because (re)creating a
# hyperlink with a TextToDisplay also (re)writes the cell's text.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:CorrectEmail@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "CorrectEmail@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:CorrectEmail@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "CorrectEmail@gmail.com")

# Update the two e-mail addresses used in the test data (E2, E3) to the new
# address, while the hyperlinks created above keep showing/linking to the
# original address.
$ws.Range("E2").Value = "CorrectEmail2@gmail.com"
$ws.Range("E3").Value = "CorrectEmail2@gmail.com"

# Update the "Actual Result" text for TC_02 (H3)
$ws.Range("H3").Value = "Error registering as the email already used"

# Move the active selection from I3 to E3
$ws.Range("E3").Select()
